$wb = $excel.ActiveWorkbook

# =======================================================================
# Sheet 1: "Measures" - insert a "Measure Folder" column before the
# existing "Measure Description" column, and refresh the sample row.
# =======================================================================
$ws1 = $wb.Worksheets.Item(1)
$tbl1 = $ws1.ListObjects.Item(1)

# ListColumns.Add() only appends at the end of the table, so grow the
# table by one column (A1:D2 -> A1:E2) and then shuffle the existing
# "Measure Description" column's header/value into the new last slot.
$tbl1.ListColumns.Add() | Out-Null

# Give the brand-new column E the same look (left/top aligned, wrapped)
# as the rest of the header row before populating it.
$ws1.Range("E1:E2").HorizontalAlignment = -4131
$ws1.Range("E1:E2").VerticalAlignment = -4160
$ws1.Range("E1:E2").WrapText = $true

$oldHeader = $ws1.Range("D1").Value2
$oldValue = $ws1.Range("D2").Value2
$ws1.Range("E1").Value = $oldHeader
$ws1.Range("E2").Value = $oldValue

# Column D becomes the new "Measure Folder" column.
$ws1.Range("D1").Value = "Measure Folder"
$ws1.Range("D2").Value = "No Folder Defined"

# Column E keeps the "Measure Description" header but gets refreshed text.
$ws1.Range("E2").Value = "This calculation is essentially adding up all of the derived line totals from a sales salesorderdetail table. This would allow a business to know total sales for all of the sales lines for a specific sales order."

# Reformat the DAX measure expression onto multiple lines.
$ws1.Range("B2").Value = "`nSUM (`n    'Sales SalesOrderDetail'[DerivedLineTotal]`n)"
# Writing a value with embedded line breaks auto-expands the row height;
# put it back to the sheet's normal (non-custom) height.
$ws1.Rows(2).AutoFit()

# Column widths: D keeps a "normal" 30-wide column, E takes over the
# wider 50 that used to belong to the description column.
$ws1.Columns("D").ColumnWidth = 29.17
$ws1.Columns("E").ColumnWidth = 49.17

# =======================================================================
# Sheet 2: "Source Information" - the Production Product / SalesOrderHeader
# / SalesTerritory rows go away; what's left is repurposed to describe the
# "Sales SalesOrderDetail" source with its rounding transformation.
# =======================================================================
$ws2 = $wb.Worksheets.Item(2)
$tbl2 = $ws2.ListObjects.Item(1)

$ws2.Range("A2").Value = 3
$ws2.Range("B2").Value = "Sales SalesOrderDetail"
$ws2.Range("F2").Value = "Sales_SalesOrderDetail"
$ws2.Range("H2").Value = "1. #`"Rounded Off`" = Table.TransformColumns(Sales_SalesOrderDetail,{{`"LineTotal`", each Number.Round(_, 2), type number}})`n"
$ws2.Range("I2").Value = "1. Rounding off the LineTotal column in the Sales_SalesOrderDetail table to two decimal places.`n"
$ws2.Rows(2).AutoFit()

# Drop the old rows 3-5 (Sales SalesOrderDetail / SalesOrderHeader / SalesTerritory).
$ws2.Rows(5).Delete()
$ws2.Rows(4).Delete()
$ws2.Rows(3).Delete()

# Restore the table/autofilter extent to its original row span.
$tbl2.Resize($ws2.Range("A1:I6"))
